# Update the Ireland ("IE") row values in the single data table on the page:
#   COD ABUNDANCE value:        0.449 -> 0.646
#   COD OVEREXPLOITATION value: 0.109 -> 0.406
#
# Both cells are targeted precisely (by table row/column) rather than with a
# document-wide text search, because "0.449" also appears (unchanged) in the
# "BE" row earlier in the table.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Row 7 (1-based) is the "IE" row:
#   IE | 0.449 | 0.109 | 0.504 | 0.214 | - | - | - | -
$row = 7

# Rebuild each cell's range via $d.Range(start, end) using the cell's own
# Start/End offsets; Find.Execute scoped this way reliably respects the
# given boundaries (searching/replacing only within that cell).

$cell1 = $table.Cell($row, 2)
$c1 = $cell1.Range
$range1 = $d.Range($c1.Start, $c1.End)
$range1.Find.Execute("0.449", $false, $false, $false, $false, $false, $true, 1, $false, "0.646", 1) | Out-Null

$cell2 = $table.Cell($row, 3)
$c2 = $cell2.Range
$range2 = $d.Range($c2.Start, $c2.End)
$range2.Find.Execute("0.109", $false, $false, $false, $false, $false, $true, 1, $false, "0.406", 1) | Out-Null
